$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert slope column (F) formulas (=1/K#) into hardcoded decimal values
# and update values to the newly re-fit standard curve numbers.
$ws.Range("F2").Value = 0.8118
$ws.Range("F3").Value = 1.406
$ws.Range("F4").Value = 1.123
$ws.Range("F7").Value = 0.26
$ws.Range("F8").Value = 0.515
$ws.Range("F9").Value = 0.338

# Update intercept column (G) values for the curves that now include an intercept
$ws.Range("G3").Value = -707.27
$ws.Range("G8").Value = -1340.12

# Update R2 value for row 3
$ws.Range("I3").Value = 0.96

# Clear the "NA" placeholders in row 5 (intercept/R2 not applicable for that slopetype)
$ws.Range("G5").ClearContents()
$ws.Range("I5").ClearContents()

# Remove the now-unused helper column K (slope_rfu_over_ugL) that the slope
# formulas used to reference
$ws.Range("K1:K11").EntireColumn.Delete()

# Update the selected cell to reflect where the editor left off
$ws.Range("G22").Select()
